$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 ("Experimental" property) -> Value column (B7) goes from blank to the
# text string "true". A leading apostrophe forces Excel to store it as text
# instead of auto-coercing it to a real Boolean; PasteSpecial(xlPasteFormats)
# from the neighbouring data cell then restores the normal data-row style
# (clearing the quote-prefix formatting the apostrophe entry would otherwise
# leave behind) so the cell keeps its original style id.
$ws.Cells.Item(7, 2).Value = "'true"
$ws.Cells.Item(6, 2).Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4122)

# Row 8 ("Date" property) -> Value column (B8) is refreshed with the new
# generation timestamp.
$ws.Cells.Item(8, 2).Value = "2025-07-21T12:46:15+00:00"

$excel.CutCopyMode = 0
